# Applies the "Add files via upload" commit:
#  - rename the worksheet tab to reflect the new extract timestamp
#  - bump the reference date (column G) for every data row from 45496 (2024-07-23)
#    to 45497 (2024-07-24)
#  - correct the one data discrepancy in row 119 (E119/H119: 617.75 -> 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to match refreshed export file name
$ws.Name = "IClientBalance-20240724-103115-"

# Determine the data extent (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 274 }

# Update the "Dt. Referencia" column (G) for every data row: 45496 -> 45497
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq 45496) {
        $cell.Value = 45497
    }
}

# Row 119 data correction: Saldo Previsto (E) and Vl. Total (H) go from 617.75 to 0
$ws.Range("E119").Value = 0
$ws.Range("H119").Value = 0
